$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.632
$ws.Range("D4").Value = -7.769
$ws.Range("A11").Value = -21.705
$ws.Range("A12").Value = -21.657
$ws.Range("D14").Value = -7.391
$ws.Range("A15").Value = -21.932
$ws.Range("D26").Value = -8.114999999999998
$ws.Range("A27").Value = -21.427
$ws.Range("A28").Value = -21.777
$ws.Range("A31").Value = -21.709
$ws.Range("D31").Value = -8.358000000000001
$ws.Range("A32").Value = -21.664
$ws.Range("D35").Value = -7.858
$ws.Range("A36").Value = -20.339
$ws.Range("D37").Value = -7.778
$ws.Range("A38").Value = -19.741
$ws.Range("D39").Value = -7.185999999999998
$ws.Range("D40").Value = -7.858
$ws.Range("D45").Value = -7.539
$ws.Range("A46").Value = -21.694
$ws.Range("D52").Value = -7.87
$ws.Range("A54").Value = -21.703
$ws.Range("A55").Value = -22.196
$ws.Range("A56").Value = -22.086
$ws.Range("D57").Value = -8.289999999999999
$ws.Range("A67").Value = -21.565
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.445
$ws.Range("A73").Value = -19.921
$ws.Range("D81").Value = -7.309
$ws.Range("A83").Value = -21.774
$ws.Range("D83").Value = -8.289
$ws.Range("A86").Value = -22.035
$ws.Range("A91").Value = -21.564
$ws.Range("A93").Value = -21.259
$ws.Range("A99").Value = -20.326
$ws.Range("D100").Value = -8.238000000000001
$ws.Range("D102").Value = -7.752
